$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value = 4427.899318504729
$ws.Range("D2").Value = 0.2542192326779488
$ws.Range("E2").Value = 3379.182024120899

$ws.Range("C3").Value = 5222.599864880444
$ws.Range("D3").Value = 0.3519302984361281
$ws.Range("E3").Value = 4677.995939524102

$ws.Range("C4").Value = 6011.77343279155
$ws.Range("D4").Value = 0.4280264039901766
$ws.Range("E4").Value = 5689.495302823623

$ws.Range("C5").Value = 6787.919122893217
$ws.Range("D5").Value = 0.4866287043590844
$ws.Range("E5").Value = 6468.460127363383

$ws.Range("C6").Value = 7573.934417395973
$ws.Range("D6").Value = 0.5627076843507866
$ws.Range("E6").Value = 7479.731850953641

$ws.Range("C7").Value = 8363.107985307079
$ws.Range("D7").Value = 0.6522993774842355
$ws.Range("E7").Value = 8670.619872117148
